# Scheduled runner update: refresh market-price-derived Leve profit columns (H-N)
# across all job sheets, per latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1785.8334
$ws.Range("I62").Value = 1785.8334
$ws.Range("K62").Value = 1785.8334
$ws.Range("M62").Value = -1161.8334
$ws.Range("H65").Value = 1785.8334
$ws.Range("I65").Value = 1785.8334
$ws.Range("K65").Value = 8929.166999999999
$ws.Range("M65").Value = -5809.166999999999
$ws.Range("H129").Value = 939.93335
$ws.Range("J129").Value = 952.0685
$ws.Range("L129").Value = 2856.2055
$ws.Range("N129").Value = -12856.2055
$ws.Range("H138").Value = 3937.653
$ws.Range("I138").Value = 851.03705
$ws.Range("J138").Value = 7725.773
$ws.Range("K138").Value = 2553.11115
$ws.Range("L138").Value = 23177.319
$ws.Range("M138").Value = 2586.88885
$ws.Range("N138").Value = -33457.319

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10681.363
$ws.Range("I45").Value = 50706
$ws.Range("J45").Value = 1787
$ws.Range("K45").Value = 50706
$ws.Range("L45").Value = 1787
$ws.Range("M45").Value = -50329
$ws.Range("N45").Value = -2541
$ws.Range("H122").Value = 1283528.8
$ws.Range("I122").Value = 1710870.8
$ws.Range("J122").Value = 1502.8
$ws.Range("K122").Value = 5132612.4
$ws.Range("L122").Value = 4508.4
$ws.Range("M122").Value = -5130162.4
$ws.Range("N122").Value = -9408.4
$ws.Range("H123").Value = 30418
$ws.Range("J123").Value = 30418
$ws.Range("L123").Value = 30418
$ws.Range("N123").Value = -40218

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H106").Value = 65503.25
$ws.Range("I106").Value = 50000
$ws.Range("J106").Value = 70671
$ws.Range("K106").Value = 50000
$ws.Range("L106").Value = 70671
$ws.Range("M106").Value = -48738
$ws.Range("N106").Value = -73195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 919.8889
$ws.Range("J16").Value = 1758
$ws.Range("L16").Value = 1758
$ws.Range("N16").Value = -2332
$ws.Range("H107").Value = 606.8333
$ws.Range("I107").Value = 497
$ws.Range("J107").Value = 716.6667
$ws.Range("K107").Value = 497
$ws.Range("L107").Value = 716.6667
$ws.Range("M107").Value = 1423
$ws.Range("N107").Value = -4556.6667
$ws.Range("H113").Value = 919.8889
$ws.Range("J113").Value = 1758
$ws.Range("L113").Value = 1758
$ws.Range("N113").Value = -6098

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 4400
$ws.Range("I70").Value = 3450
$ws.Range("J70").Value = 4875
$ws.Range("K70").Value = 10350
$ws.Range("L70").Value = 14625
$ws.Range("M70").Value = -10035
$ws.Range("N70").Value = -15255
$ws.Range("H73").Value = 4400
$ws.Range("I73").Value = 3450
$ws.Range("J73").Value = 4875
$ws.Range("K73").Value = 10350
$ws.Range("L73").Value = 14625
$ws.Range("M73").Value = -9258
$ws.Range("N73").Value = -16809
$ws.Range("H94").Value = 3200
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 4550
$ws.Range("K94").Value = 1500
$ws.Range("L94").Value = 13650
$ws.Range("M94").Value = -824
$ws.Range("N94").Value = -15002
$ws.Range("H112").Value = 2403.9
$ws.Range("I112").Value = 675.6667
$ws.Range("J112").Value = 2595.926
$ws.Range("K112").Value = 2027.0001
$ws.Range("L112").Value = 7787.778
$ws.Range("M112").Value = -919.0001
$ws.Range("N112").Value = -10003.778
$ws.Range("H113").Value = 5008839
$ws.Range("I113").Value = 10000474
$ws.Range("J113").Value = 1443385.1
$ws.Range("K113").Value = 30001422
$ws.Range("L113").Value = 4330155.300000001
$ws.Range("M113").Value = -29999252
$ws.Range("N113").Value = -4334495.300000001
$ws.Range("H123").Value = 6751.778
$ws.Range("I123").Value = 2933.3333
$ws.Range("J123").Value = 8661
$ws.Range("K123").Value = 8799.999899999999
$ws.Range("L123").Value = 25983
$ws.Range("M123").Value = -6349.999899999999
$ws.Range("N123").Value = -30883
$ws.Range("H134").Value = 9160.205
$ws.Range("I134").Value = 7816.533
$ws.Range("K134").Value = 23449.599
$ws.Range("M134").Value = -18379.599

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 125001190
$ws.Range("I113").Value = 250000800
$ws.Range("J113").Value = 1578.25
$ws.Range("K113").Value = 250000800
$ws.Range("L113").Value = 1578.25
$ws.Range("M113").Value = -249998630
$ws.Range("N113").Value = -5918.25
$ws.Range("H122").Value = 4254325
$ws.Range("I122").Value = 5402926
$ws.Range("J122").Value = 3335443.8
$ws.Range("K122").Value = 16208778
$ws.Range("L122").Value = 10006331.4
$ws.Range("M122").Value = -16206328
$ws.Range("N122").Value = -10011231.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 742.3333
$ws.Range("I16").Value = 761.125
$ws.Range("J16").Value = 592
$ws.Range("K16").Value = 761.125
$ws.Range("L16").Value = 592
$ws.Range("M16").Value = -591.125
$ws.Range("N16").Value = -932
$ws.Range("H61").Value = 2411.111
$ws.Range("I61").Value = 2462.5
$ws.Range("K61").Value = 2462.5
$ws.Range("M61").Value = -2260.5
$ws.Range("H113").Value = 2411.111
$ws.Range("I113").Value = 2462.5
$ws.Range("K113").Value = 2462.5
$ws.Range("M113").Value = -292.5
$ws.Range("H122").Value = 5430705
$ws.Range("I122").Value = 5954822.5
$ws.Range("J122").Value = 3334235
$ws.Range("K122").Value = 17864467.5
$ws.Range("L122").Value = 10002705
$ws.Range("M122").Value = -17862017.5
$ws.Range("N122").Value = -10007605
$ws.Range("H136").Value = 3811
$ws.Range("I136").Value = 2676.3462
$ws.Range("J136").Value = 7088.8887
$ws.Range("K136").Value = 8029.0386
$ws.Range("L136").Value = 21266.6661
$ws.Range("M136").Value = -5479.0386
$ws.Range("N136").Value = -26366.6661

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 250000620
$ws.Range("I107").Value = 333333900
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1000001700
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = -999999780
$ws.Range("N107").Value = -6240
$ws.Range("H122").Value = 1979.091
$ws.Range("I122").Value = 1971.25
$ws.Range("K122").Value = 5913.75
$ws.Range("M122").Value = -3463.75
$ws.Range("H123").Value = 29840.908
$ws.Range("J123").Value = 29840.908
$ws.Range("L123").Value = 29840.908
$ws.Range("N123").Value = -39640.908
$ws.Range("H132").Value = 2170.5833
$ws.Range("I132").Value = 1459.4667
$ws.Range("J132").Value = 3355.7778
$ws.Range("K132").Value = 4378.4001
$ws.Range("L132").Value = 10067.3334
$ws.Range("M132").Value = -1848.4001
$ws.Range("N132").Value = -15127.3334
